$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gmail")

# Add a new row of test data (Nisha / 123456) below the existing Gmail entry
$ws.Range("A3").Value = "Nisha"
$ws.Range("B3").Value = 123456

# Move the active selection, matching where the author last left the cursor
$ws.Range("B10").Select() | Out-Null
